# Add a new "Save" column (H) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: text "Save", formatted like the neighboring header G1
# (bold font, thin border, centered alignment) by copying G1's formatting.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells H2:H6 - plain numeric 0/1 flags, no special style (same as
# the other numeric data columns B:G).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
